$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "% of Q Drop's" column header, right after the existing
# "% of F's" column (H).
$ws.Cells.Item(1, 9).Value = "% of Q Drop's"

# Populate the new column with "0.00%" text values for each row that
# already holds grade-distribution figures. A leading apostrophe forces
# the value to be stored as text (matching the existing percentage
# columns) rather than being reinterpreted as a number.
$dataRows = @(3, 6, 9, 12, 13, 16)
foreach ($r in $dataRows) {
    $ws.Cells.Item($r, 9).Value = "'0.00%"
}

# Clear the automatically-inferred "text/quote-prefix" formatting so the
# new cells keep the same (default) style as the rest of the sheet.
foreach ($r in $dataRows) {
    $ws.Cells.Item($r, 9).Style = "Normal"
}
